# Delete row 731 (the llama post) and shift all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(731).Delete()
